$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4275.108
$ws.Range("I11").Value = 4275.108
$ws.Range("K11").Value = 4275.108
$ws.Range("M11").Value = -4135.108
$ws.Range("H12").Value = 360.66666
$ws.Range("K12").Value = 425
$ws.Range("M12").Value = -255
$ws.Range("I12").Value = 425
$ws.Range("L12").Value = 39
$ws.Range("J12").Value = 39
$ws.Range("N12").Value = -379
$ws.Range("M15").Value = -5796.0524
$ws.Range("I15").Value = 1988.3508
$ws.Range("H15").Value = 1988.3508
$ws.Range("K15").Value = 5965.0524
$ws.Range("H34").Value = 8005.5
$ws.Range("K34").Value = 8005.5
$ws.Range("M34").Value = -7802.5
$ws.Range("I34").Value = 8005.5
$ws.Range("I36").Value = 8005.5
$ws.Range("M36").Value = -7290.5
$ws.Range("K36").Value = 8005.5
$ws.Range("H36").Value = 8005.5
$ws.Range("H47").Value = 67724.75
$ws.Range("M47").Value = -46478
$ws.Range("K47").Value = 47450
$ws.Range("I47").Value = 47450
$ws.Range("H54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("M54").Value = -4514
$ws.Range("K62").Value = 53339932
$ws.Range("H62").Value = 53339932
$ws.Range("M62").Value = -53339308
$ws.Range("I62").Value = 53339932
$ws.Range("H65").Value = 53339932
$ws.Range("K65").Value = 266699660
$ws.Range("I65").Value = 53339932
$ws.Range("M65").Value = -266696540
$ws.Range("H76").Value = 10004425
$ws.Range("I76").Value = 16670333
$ws.Range("K76").Value = 16670333
$ws.Range("L76").Value = 5562
$ws.Range("N76").Value = -6192
$ws.Range("J76").Value = 5562
$ws.Range("M76").Value = -16670018
$ws.Range("I79").Value = 16670333
$ws.Range("K79").Value = 16670333
$ws.Range("J79").Value = 5562
$ws.Range("L79").Value = 5562
$ws.Range("M79").Value = -16669241
$ws.Range("H79").Value = 10004425
$ws.Range("N79").Value = -7746
$ws.Range("K80").Value = 1126.875
$ws.Range("H80").Value = 1238.9166
$ws.Range("N80").Value = -7007.6875
$ws.Range("I80").Value = 375.625
$ws.Range("L80").Value = 5011.6875
$ws.Range("M80").Value = -128.875
$ws.Range("J80").Value = 1670.5625
$ws.Range("M83").Value = 1611.375
$ws.Range("J83").Value = 1670.5625
$ws.Range("L83").Value = 15035.0625
$ws.Range("H83").Value = 1238.9166
$ws.Range("I83").Value = 375.625
$ws.Range("K83").Value = 3380.625
$ws.Range("N83").Value = -25019.0625
$ws.Range("I88").Value = 3003
$ws.Range("K88").Value = 3003
$ws.Range("H88").Value = 6501.5
$ws.Range("M88").Value = -2597
$ws.Range("H91").Value = 6501.5
$ws.Range("I91").Value = 3003
$ws.Range("K91").Value = 3003
$ws.Range("M91").Value = -1599
$ws.Range("K106").Value = 3350
$ws.Range("I106").Value = 3350
$ws.Range("M106").Value = -2719
$ws.Range("H106").Value = 3416.6667
$ws.Range("L113").Value = 202041
$ws.Range("H113").Value = 113500.11
$ws.Range("I113").Value = 2824
$ws.Range("M113").Value = 430
$ws.Range("K113").Value = 2824
$ws.Range("N113").Value = -208549
$ws.Range("J113").Value = 202041
$ws.Range("H133").Value = 119949.5
$ws.Range("J133").Value = 119949.5
$ws.Range("N133").Value = -130069.5
$ws.Range("L133").Value = 119949.5
$ws.Range("H135").Value = 2601.0715
$ws.Range("I135").Value = 2207
$ws.Range("M135").Value = -17328
$ws.Range("L135").Value = 28138.5
$ws.Range("K135").Value = 19863
$ws.Range("N135").Value = -33208.5
$ws.Range("J135").Value = 3126.5
$ws.Range("N137").Value = -16647
$ws.Range("I137").Value = 4123.75
$ws.Range("L137").Value = 11547
$ws.Range("J137").Value = 3849
$ws.Range("K137").Value = 12371.25
$ws.Range("M137").Value = -9821.25
$ws.Range("H137").Value = 3958.9
$ws.Range("M138").Value = -2270
$ws.Range("I138").Value = 2470
$ws.Range("H138").Value = 2558.762
$ws.Range("K138").Value = 7410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6280
$ws.Range("M3").Value = -4021
$ws.Range("K3").Value = 4136
$ws.Range("I3").Value = 4136
$ws.Range("L8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H12").Value = 44058.168
$ws.Range("N12").Value = -66097
$ws.Range("K12").Value = 672.5
$ws.Range("M12").Value = -499.5
$ws.Range("I12").Value = 672.5
$ws.Range("L12").Value = 65751
$ws.Range("J12").Value = 65751
$ws.Range("K25").Value = 501
$ws.Range("H25").Value = 1045.5454
$ws.Range("I25").Value = 501
$ws.Range("M25").Value = -99
$ws.Range("M61").Value = -7421.615
$ws.Range("L61").Value = 16083.333
$ws.Range("J61").Value = 16083.333
$ws.Range("K61").Value = 7633.615
$ws.Range("H61").Value = 8507.725
$ws.Range("I61").Value = 7633.615
$ws.Range("N61").Value = -16507.333
$ws.Range("N63").Value = -11104.368
$ws.Range("K63").Value = 2442.75
$ws.Range("J63").Value = 9732.368
$ws.Range("I63").Value = 2442.75
$ws.Range("H63").Value = 6399.971
$ws.Range("M63").Value = -1756.75
$ws.Range("L63").Value = 9732.368
$ws.Range("M66").Value = -8781.75
$ws.Range("J66").Value = 9732.368
$ws.Range("H66").Value = 6399.971
$ws.Range("N66").Value = -55525.84
$ws.Range("I66").Value = 2442.75
$ws.Range("L66").Value = 48661.84
$ws.Range("K66").Value = 12213.75
$ws.Range("I74").Value = 2001662.2
$ws.Range("M74").Value = -2000788.2
$ws.Range("H74").Value = 913420.8
$ws.Range("L74").Value = 6553
$ws.Range("J74").Value = 6553
$ws.Range("K74").Value = 2001662.2
$ws.Range("N74").Value = -8301
$ws.Range("N77").Value = -41501
$ws.Range("K77").Value = 10008311
$ws.Range("H77").Value = 913420.8
$ws.Range("M77").Value = -10003943
$ws.Range("I77").Value = 2001662.2
$ws.Range("L77").Value = 32765
$ws.Range("J77").Value = 6553
$ws.Range("N88").Value = -2218.2
$ws.Range("H88").Value = 1458.4166
$ws.Range("L88").Value = 1406.2
$ws.Range("J88").Value = 1406.2
$ws.Range("H91").Value = 1458.4166
$ws.Range("L91").Value = 1406.2
$ws.Range("N91").Value = -4214.2
$ws.Range("J91").Value = 1406.2
$ws.Range("J101").Value = 58249.5
$ws.Range("H101").Value = 58249.5
$ws.Range("N101").Value = -64739.5
$ws.Range("L101").Value = 58249.5
$ws.Range("N102").Value = -83671584
$ws.Range("K102").Value = 4239.6665
$ws.Range("H102").Value = 41836290
$ws.Range("I102").Value = 4239.6665
$ws.Range("M102").Value = -2617.6665
$ws.Range("L102").Value = 83668340
$ws.Range("J102").Value = 83668340
$ws.Range("I110").Value = 15625765
$ws.Range("H110").Value = 11365329
$ws.Range("K110").Value = 15625765
$ws.Range("M110").Value = -15623720
$ws.Range("K122").Value = 8704.5
$ws.Range("H122").Value = 3854.9092
$ws.Range("M122").Value = -6254.5
$ws.Range("I122").Value = 2901.5
$ws.Range("H136").Value = 8507.725
$ws.Range("K136").Value = 22900.845
$ws.Range("M136").Value = -20350.845
$ws.Range("N136").Value = -53349.999
$ws.Range("I136").Value = 7633.615
$ws.Range("L136").Value = 48249.999
$ws.Range("J136").Value = 16083.333
$ws.Range("H138").Value = 119999
$ws.Range("J138").Value = 119999
$ws.Range("N138").Value = -130279
$ws.Range("L138").Value = 119999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1452.2693
$ws.Range("K94").Value = 1276.5294
$ws.Range("M94").Value = -825.5293999999999
$ws.Range("N94").Value = -2686.2222
$ws.Range("I94").Value = 1276.5294
$ws.Range("J94").Value = 1784.2222
$ws.Range("L94").Value = 1784.2222
$ws.Range("H107").Value = 1323.45
$ws.Range("L107").Value = 1358.6
$ws.Range("J107").Value = 1358.6
$ws.Range("I107").Value = 1288.3
$ws.Range("N107").Value = -5198.6
$ws.Range("K107").Value = 1288.3
$ws.Range("M107").Value = 631.7
$ws.Range("J134").Value = 8728.111000000001
$ws.Range("K134").Value = 9719.625
$ws.Range("N134").Value = -31254.333
$ws.Range("I134").Value = 3239.875
$ws.Range("L134").Value = 26184.333
$ws.Range("H134").Value = 4247.9185
$ws.Range("M134").Value = -7184.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1999.6666
$ws.Range("M4").Value = -388
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("N21").Value = -10466.75
$ws.Range("K21").Value = 2001
$ws.Range("I21").Value = 2001
$ws.Range("H21").Value = 8397.6
$ws.Range("J21").Value = 9996.75
$ws.Range("L21").Value = 9996.75
$ws.Range("M21").Value = -1766
$ws.Range("H31").Value = 66672736
$ws.Range("I31").Value = 200001840
$ws.Range("M31").Value = -200001545
$ws.Range("K31").Value = 200001840
$ws.Range("H34").Value = 66672736
$ws.Range("K34").Value = 200001840
$ws.Range("M34").Value = -200001638
$ws.Range("I34").Value = 200001840
$ws.Range("K62").Value = 4250
$ws.Range("H62").Value = 8815.143
$ws.Range("M62").Value = -3626
$ws.Range("I62").Value = 4250
$ws.Range("H65").Value = 8815.143
$ws.Range("K65").Value = 21250
$ws.Range("I65").Value = 4250
$ws.Range("M65").Value = -18130
$ws.Range("H107").Value = 1250
$ws.Range("I107").Value = 1250
$ws.Range("K107").Value = 1250
$ws.Range("M107").Value = 670
$ws.Range("N109").Value = -41222.312
$ws.Range("J109").Value = 39142.312
$ws.Range("L109").Value = 39142.312
$ws.Range("H109").Value = 39142.312
$ws.Range("K122").Value = 1004881.02
$ws.Range("H122").Value = 145697.28
$ws.Range("M122").Value = -1002431.02
$ws.Range("I122").Value = 334960.34
$ws.Range("J134").Value = 10275
$ws.Range("K134").Value = 8485.200000000001
$ws.Range("N134").Value = -35895
$ws.Range("I134").Value = 2828.4
$ws.Range("L134").Value = 30825
$ws.Range("H134").Value = 3704.4707
$ws.Range("M134").Value = -5950.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2400985.5
$ws.Range("M4").Value = -8400841.399999999
$ws.Range("I4").Value = 2800317.8
$ws.Range("K4").Value = 8400953.399999999
$ws.Range("L4").Value = 14976
$ws.Range("N4").Value = -15200
$ws.Range("J4").Value = 4992
$ws.Range("J9").Value = 4233
$ws.Range("N9").Value = -13147
$ws.Range("H9").Value = 4233
$ws.Range("L9").Value = 12699
$ws.Range("K23").Value = 375
$ws.Range("M23").Value = -140
$ws.Range("N23").Value = -1027.14284
$ws.Range("J23").Value = 185.71428
$ws.Range("H23").Value = 178.125
$ws.Range("L23").Value = 557.14284
$ws.Range("I23").Value = 125
$ws.Range("J37").Value = 76604.414
$ws.Range("L37").Value = 229813.242
$ws.Range("H37").Value = 76604.414
$ws.Range("N37").Value = -230037.242
$ws.Range("N39").Value = -17732.1432
$ws.Range("J39").Value = 5714.7144
$ws.Range("H39").Value = 4505.8887
$ws.Range("L39").Value = 17144.1432
$ws.Range("H43").Value = 5334.6665
$ws.Range("L43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H56").Value = 11583.667
$ws.Range("M56").Value = -11053.667
$ws.Range("I56").Value = 11583.667
$ws.Range("K56").Value = 11583.667
$ws.Range("I60").Value = 199.6
$ws.Range("H60").Value = 266.33334
$ws.Range("K60").Value = 598.8
$ws.Range("M60").Value = -347.8
$ws.Range("L113").Value = 5608.799999999999
$ws.Range("H113").Value = 1528.1428
$ws.Range("I113").Value = 674.5
$ws.Range("M113").Value = 146.5
$ws.Range("K113").Value = 2023.5
$ws.Range("N113").Value = -9948.799999999999
$ws.Range("J113").Value = 1869.6
$ws.Range("J122").Value = 1000
$ws.Range("N122").Value = -13900
$ws.Range("H122").Value = 1100
$ws.Range("L122").Value = 9000
$ws.Range("I129").Value = 485.85715
$ws.Range("H129").Value = 45455496
$ws.Range("M129").Value = 3542.42855
$ws.Range("K129").Value = 1457.57145
$ws.Range("H132").Value = 49563.81
$ws.Range("M132").Value = -609824.6
$ws.Range("K132").Value = 612354.6
$ws.Range("I132").Value = 68039.39999999999
$ws.Range("H136").Value = 4890.75
$ws.Range("N136").Value = -29749.5
$ws.Range("L136").Value = 19549.5
$ws.Range("J136").Value = 6516.5
$ws.Range("I138").Value = 0
$ws.Range("H138").Value = 2505274.8
$ws.Range("K138").Value = 0
$ws.Range("J138").Value = 2505274.8
$ws.Range("N138").Value = -7526104.399999999
$ws.Range("L138").Value = 7515824.399999999
$ws.Range("M138").ClearContents()
$ws.Range("I140").Value = 982.5
$ws.Range("K140").Value = 2947.5
$ws.Range("N140").Value = -15092.5
$ws.Range("J140").Value = 1577.5
$ws.Range("H140").Value = 1280
$ws.Range("M140").Value = 2232.5
$ws.Range("L140").Value = 4732.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 9899.5
$ws.Range("J38").Value = 9899.5
$ws.Range("L38").Value = 9899.5
$ws.Range("N38").Value = -10825.5
$ws.Range("K80").Value = 4666
$ws.Range("H80").Value = 6309.6665
$ws.Range("I80").Value = 4666
$ws.Range("M80").Value = -3668
$ws.Range("M83").Value = -18338
$ws.Range("H83").Value = 6309.6665
$ws.Range("I83").Value = 4666
$ws.Range("K83").Value = 23330
$ws.Range("M97").Value = -1164.6666
$ws.Range("H97").Value = 1942.125
$ws.Range("K97").Value = 1660.6666
$ws.Range("I97").Value = 1660.6666
$ws.Range("J122").Value = 49999
$ws.Range("N122").Value = -154897
$ws.Range("K122").Value = 2664
$ws.Range("H122").Value = 25443.5
$ws.Range("L122").Value = 149997
$ws.Range("I122").Value = 888
$ws.Range("M122").Value = -214
$ws.Range("I126").Value = 1999.7858
$ws.Range("H126").Value = 3190.7144
$ws.Range("M126").Value = -3529.357400000001
$ws.Range("K126").Value = 5999.357400000001
$ws.Range("H132").Value = 6249.923
$ws.Range("M132").Value = -10765.319
$ws.Range("K132").Value = 13295.319
$ws.Range("I132").Value = 4431.773

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L22").Value = 4719.579
$ws.Range("J22").Value = 4719.579
$ws.Range("N22").Value = -5309.579
$ws.Range("H22").Value = 3571.162
$ws.Range("J27").Value = 4719.579
$ws.Range("L27").Value = 4719.579
$ws.Range("H27").Value = 3571.162
$ws.Range("N27").Value = -4933.579
$ws.Range("H46").Value = 6903.1377
$ws.Range("M46").Value = -1596.2858
$ws.Range("I46").Value = 1784.2858
$ws.Range("K46").Value = 1784.2858
$ws.Range("L61").Value = 2651
$ws.Range("J61").Value = 2651
$ws.Range("H61").Value = 2112
$ws.Range("N61").Value = -3055
$ws.Range("N93").Value = -7621.5
$ws.Range("J93").Value = 5125.5
$ws.Range("M93").Value = -3415
$ws.Range("I93").Value = 4663
$ws.Range("H93").Value = 4805.3076
$ws.Range("L93").Value = 5125.5
$ws.Range("K93").Value = 4663
$ws.Range("L113").Value = 2651
$ws.Range("H113").Value = 2112
$ws.Range("N113").Value = -6991
$ws.Range("J113").Value = 2651
$ws.Range("H132").Value = 2730.457
$ws.Range("M132").Value = -3329
$ws.Range("J132").Value = 4674.1
$ws.Range("K132").Value = 5859
$ws.Range("N132").Value = -19082.3
$ws.Range("L132").Value = 14022.3
$ws.Range("I132").Value = 1953
$ws.Range("N140").Value = -119860
$ws.Range("J140").Value = 109500
$ws.Range("H140").Value = 81333.336
$ws.Range("L140").Value = 109500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N18").Value = -34096
$ws.Range("J18").Value = 33750
$ws.Range("H18").Value = 32222.223
$ws.Range("L18").Value = 33750
$ws.Range("N109").Value = -63396.5
$ws.Range("J109").Value = 60622.5
$ws.Range("L109").Value = 60622.5
$ws.Range("H109").Value = 52098
$ws.Range("K122").Value = 15900
$ws.Range("H122").Value = 4971.4287
$ws.Range("M122").Value = -13450
$ws.Range("I122").Value = 5300
$ws.Range("H132").Value = 8234.157999999999
$ws.Range("M132").Value = -17996
$ws.Range("K132").Value = 20526
$ws.Range("I132").Value = 6842
